$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Shift the "cube ordinate=" (G4/H4) and "calc r=" (G5/H5) label+formula
#    pairs one column to the right, to H4/I4 and H5/I5 respectively.
#    Remember the formulas first (Copy() only duplicates cached values, not
#    formulas, so the numeric formulas must be moved by hand).
# ---------------------------------------------------------------------------
$h4Formula = $ws.Range("H4").Formula
$h5Formula = $ws.Range("H5").Formula

$ws.Range("G4").Copy($ws.Range("H4"))
$ws.Range("G5").Copy($ws.Range("H5"))
$ws.Range("G4").Clear()
$ws.Range("G5").Clear()

$ws.Range("I4").Formula = $h4Formula
$ws.Range("I5").Formula = $h5Formula

# ---------------------------------------------------------------------------
# 2. Defined name "ordScale" now refers to I4 instead of H4. Re-creating the
#    name (rather than only changing RefersTo) makes existing/new formulas
#    that use it recalculate against the new target cell.
# ---------------------------------------------------------------------------
$wb.Names.Item("ordScale").Delete()
$wb.Names.Add("ordScale", "=Sheet1!`$I`$4")

# ---------------------------------------------------------------------------
# 3. New "volume =" / "surface area=" rows (6 & 7), using the same label
#    style as the other right-aligned headers in column G/H (H5, moved
#    above, already carries that style).
# ---------------------------------------------------------------------------
$ws.Range("H5").Copy($ws.Range("H6"))
$ws.Range("H6").Value = "volume ="
$ws.Range("I6").Formula = "=(2*z*ordScale)^3 *(15*7*SQRT(5))/4"
$ws.Range("J6").Formula = "=I6/12"
$ws.Range("K6").Formula = "=J6/5"

$ws.Range("H5").Copy($ws.Range("H7"))
$ws.Range("H7").Value = "surface area="
$ws.Range("I7").Formula = "=3*(2*z*ordScale)^2*SQRT(25+10*SQRT(5))"
$ws.Range("J7").Formula = "=I7/12"

# ---------------------------------------------------------------------------
# 4. "Center:" average blocks for the four tetrahedra face groups, added at
#    rows 17, 25 and 33.
# ---------------------------------------------------------------------------
function Add-CenterRow($row, $srcFirst, $srcLast) {
    foreach ($col0 in @("I", "P", "W", "AD")) {
        $ws.Range("$col0$row").Value = "Center:"
    }

    $groups = @(
        @("J", "K", "L", "M"),
        @("Q", "R", "S", "T"),
        @("X", "Y", "Z", "AA"),
        @("AE", "AF", "AG", "AH")
    )
    foreach ($cols in $groups) {
        foreach ($col in $cols) {
            $ws.Range("$col$row").Formula = "=AVERAGE($col$srcFirst`:$col$srcLast)"
        }
    }
}

Add-CenterRow 17 12 16
Add-CenterRow 25 20 24
Add-CenterRow 33 28 32

# ---------------------------------------------------------------------------
# 5. Sheet view: scroll/selection state
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$ws.Range("W33:AA33").Select()
$win.ScrollRow = 13
$win.ScrollColumn = 12
